$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the recalculated time-study percentages/coefficients (Тпз block)
$ws.Range("C23").Value = "0.83"
$ws.Range("C24").Value = "82.5"
$ws.Range("C25").Value = "0.83"
$ws.Range("C26").Value = "2.08"
$ws.Range("C28").Value = "1.25"

# Update the recalculated total (Тобс/ИТОГО row)
$ws.Range("B29").Value = 480

# Move the active selection (mirrors the saved cursor position in the workbook)
$ws.Range("A35").Select() | Out-Null
